$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 211.2
$ws.Range("I33").Value = 195.27272
$ws.Range("K33").Value = 195.27272
$ws.Range("M33").Value = 33.72728000000001
$ws.Range("H40").Value = 5180411
$ws.Range("I40").Value = 6251795
$ws.Range("J40").Value = 2501951
$ws.Range("K40").Value = 6251795
$ws.Range("L40").Value = 2501951
$ws.Range("M40").Value = -6251620
$ws.Range("N40").Value = -2502301
$ws.Range("H70").Value = 1389.5
$ws.Range("I70").Value = 2640.6
$ws.Range("J70").Value = 694.44446
$ws.Range("K70").Value = 7921.799999999999
$ws.Range("L70").Value = 2083.33338
$ws.Range("M70").Value = -7651.799999999999
$ws.Range("N70").Value = -2623.33338
$ws.Range("H73").Value = 1389.5
$ws.Range("I73").Value = 2640.6
$ws.Range("J73").Value = 694.44446
$ws.Range("K73").Value = 7921.799999999999
$ws.Range("L73").Value = 2083.33338
$ws.Range("M73").Value = -6985.799999999999
$ws.Range("N73").Value = -3955.33338
$ws.Range("H116").Value = 2220.8333
$ws.Range("J116").Value = 2392.2307
$ws.Range("L116").Value = 2392.2307
$ws.Range("N116").Value = -9276.2307

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 11385.714
$ws.Range("I36").Value = 10175
$ws.Range("J36").Value = 13000
$ws.Range("K36").Value = 10175
$ws.Range("L36").Value = 13000
$ws.Range("M36").Value = -9829
$ws.Range("N36").Value = -13692
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("N37").ClearContents()
$ws.Range("H102").Value = 1315.7368
$ws.Range("I102").Value = 1249.9445
$ws.Range("J102").Value = 2500
$ws.Range("K102").Value = 1249.9445
$ws.Range("L102").Value = 2500
$ws.Range("M102").Value = 372.0554999999999
$ws.Range("N102").Value = -5744

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 4710
$ws.Range("I26").Value = 4710
$ws.Range("K26").Value = 4710
$ws.Range("M26").Value = -4418
$ws.Range("H28").Value = 25000
$ws.Range("J28").Value = 25000
$ws.Range("L28").Value = 25000
$ws.Range("N28").Value = -25588
$ws.Range("H122").Value = 30780
$ws.Range("J122").Value = 30780
$ws.Range("L122").Value = 30780
$ws.Range("N122").Value = -40580
$ws.Range("H123").Value = 43389.5
$ws.Range("J123").Value = 43389.5
$ws.Range("L123").Value = 43389.5
$ws.Range("N123").Value = -53189.5
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("H127").Value = 49999.5
$ws.Range("J127").Value = 49999.5
$ws.Range("L127").Value = 49999.5
$ws.Range("N127").Value = -59919.5
$ws.Range("H129").Value = 49999
$ws.Range("J129").Value = 49999
$ws.Range("L129").Value = 49999
$ws.Range("N129").Value = -59999
$ws.Range("H130").Value = 30780
$ws.Range("J130").Value = 30780
$ws.Range("L130").Value = 30780
$ws.Range("N130").Value = -40820
$ws.Range("H131").Value = 39000
$ws.Range("J131").Value = 39000
$ws.Range("L131").Value = 39000
$ws.Range("N131").Value = -49080

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 5065
$ws.Range("I19").Value = 5065
$ws.Range("K19").Value = 5065
$ws.Range("M19").Value = -4895
$ws.Range("H20").Value = 50194.25
$ws.Range("J20").Value = 50194.25
$ws.Range("L20").Value = 50194.25
$ws.Range("N20").Value = -50666.25
$ws.Range("H23").Value = 58728.75
$ws.Range("I23").Value = 3000
$ws.Range("J23").Value = 66690
$ws.Range("K23").Value = 3000
$ws.Range("L23").Value = 66690
$ws.Range("M23").Value = -2760
$ws.Range("N23").Value = -67170
$ws.Range("H24").Value = 5065
$ws.Range("I24").Value = 5065
$ws.Range("K24").Value = 5065
$ws.Range("M24").Value = -4895
$ws.Range("H25").Value = 5000
$ws.Range("I25").Value = 5000
$ws.Range("K25").Value = 5000
$ws.Range("M25").Value = -4826
$ws.Range("H27").Value = 58728.75
$ws.Range("I27").Value = 3000
$ws.Range("J27").Value = 66690
$ws.Range("K27").Value = 3000
$ws.Range("L27").Value = 66690
$ws.Range("M27").Value = -2808
$ws.Range("N27").Value = -67074
$ws.Range("H30").Value = 50194.25
$ws.Range("J30").Value = 50194.25
$ws.Range("L30").Value = 50194.25
$ws.Range("N30").Value = -50376.25
$ws.Range("H58").Value = 50000956
$ws.Range("I58").Value = 100001100
$ws.Range("J58").Value = 807
$ws.Range("K58").Value = 100001100
$ws.Range("L58").Value = 807
$ws.Range("M58").Value = -100000897
$ws.Range("N58").Value = -1213
$ws.Range("H124").Value = 15206.25
$ws.Range("J124").Value = 15206.25
$ws.Range("L124").Value = 15206.25
$ws.Range("N124").Value = -20116.25
$ws.Range("H125").Value = 10000
$ws.Range("J125").Value = 10000
$ws.Range("L125").Value = 10000
$ws.Range("N125").Value = -14920
$ws.Range("H128").Value = 50194.25
$ws.Range("J128").Value = 50194.25
$ws.Range("L128").Value = 50194.25
$ws.Range("N128").Value = -60154.25
$ws.Range("H129").Value = 49999
$ws.Range("J129").Value = 49999
$ws.Range("L129").Value = 49999
$ws.Range("N129").Value = -59999
$ws.Range("H130").Value = 30390
$ws.Range("J130").Value = 30390
$ws.Range("L130").Value = 30390
$ws.Range("N130").Value = -40430
$ws.Range("H131").Value = 36000
$ws.Range("J131").Value = 36000
$ws.Range("L131").Value = 36000
$ws.Range("N131").Value = -46080
$ws.Range("H134").Value = 1340.0476
$ws.Range("I134").Value = 1321.6428
$ws.Range("J134").Value = 1376.8572
$ws.Range("K134").Value = 3964.9284
$ws.Range("L134").Value = 4130.571599999999
$ws.Range("M134").Value = -1429.9284
$ws.Range("N134").Value = -9200.571599999999
$ws.Range("H136").Value = 50000956
$ws.Range("I136").Value = 100001100
$ws.Range("J136").Value = 807
$ws.Range("K136").Value = 300003300
$ws.Range("L136").Value = 2421
$ws.Range("M136").Value = -300000750
$ws.Range("N136").Value = -7521

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 5382
$ws.Range("I87").Value = 3509.3333
$ws.Range("K87").Value = 10527.9999
$ws.Range("M87").Value = -9279.999899999999
$ws.Range("H90").Value = 5382
$ws.Range("I90").Value = 3509.3333
$ws.Range("K90").Value = 31583.9997
$ws.Range("M90").Value = -25343.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 776.1852
$ws.Range("I97").Value = 796.75
$ws.Range("J97").Value = 717.4286
$ws.Range("K97").Value = 796.75
$ws.Range("L97").Value = 717.4286
$ws.Range("M97").Value = -300.75
$ws.Range("N97").Value = -1709.4286

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 845881.8
$ws.Range("I22").Value = 2111504.5
$ws.Range("J22").Value = 2133.3333
$ws.Range("K22").Value = 2111504.5
$ws.Range("L22").Value = 2133.3333
$ws.Range("M22").Value = -2111209.5
$ws.Range("N22").Value = -2723.3333
$ws.Range("H24").Value = 2750
$ws.Range("J24").Value = 2750
$ws.Range("L24").Value = 2750
$ws.Range("N24").Value = -3436
$ws.Range("H27").Value = 845881.8
$ws.Range("I27").Value = 2111504.5
$ws.Range("J27").Value = 2133.3333
$ws.Range("K27").Value = 2111504.5
$ws.Range("L27").Value = 2133.3333
$ws.Range("M27").Value = -2111397.5
$ws.Range("N27").Value = -2347.3333
$ws.Range("H100").Value = 1449.9166
$ws.Range("I100").Value = 1400
$ws.Range("J100").Value = 1459.9
$ws.Range("K100").Value = 1400
$ws.Range("L100").Value = 1459.9
$ws.Range("M100").Value = -859
$ws.Range("N100").Value = -2541.9
